$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "94.749.84"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.590.41"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.97%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "656.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.63%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.45"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.402"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.72%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.984"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.57%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.584.71"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.85%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.200"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.24%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.20"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.44%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.284.65"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.081.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.86%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000251"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.80%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.602.28"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.27%  "

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.90"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.47%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +10.46%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.55"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.21%  "

# Row 23
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "503.85"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.91%  "

# Row 24
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.475"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -9.46%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000194"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.54"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +10.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.780.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.41"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.14"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +15.60%  "

# Row 31
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.01%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.20"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.68%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.137"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.18%  "

# Row 34
$ws.Range("E34").Value = "  -0.51%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.175"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.70"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.87%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.554"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "570.56"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.57%  "

# Row 41
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.148"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.19%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.913"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "34.69"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +43.94%  "

# Row 45
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.63"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.42%  "

# Row 46
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.24%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.21"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0408"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.40%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.48"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.95%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.44"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.55%  "
